$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 129; this pushes the existing rows 129-207
# down to 130-208 (preserving their data, formatting and the D-column
# date style) and leaves the new row 129 blank and ready to be filled
# in with the newly reported price observation.
$ws.Rows(129).Insert()

# Populate the newly inserted row 129 with the new observation. The
# surrounding (unchanged) columns A, B, C, E, F, G, H, N, O, Q, R keep the
# same constant values used throughout this sheet, so we only need to
# set them explicitly for this brand-new row.
$ws.Range("A129").Value = 11
$ws.Range("B129").Value = "Vega Monumental Concepción"
$ws.Range("C129").Value = "Bíobío"
$ws.Range("D129").Value = 44572
$ws.Range("D129").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E129").Value = 8
$ws.Range("F129").Value = 100112023
$ws.Range("G129").Value = "Brócoli"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 270
$ws.Range("K129").Value = 850
$ws.Range("L129").Value = 900
$ws.Range("M129").Value = 878
$ws.Range("N129").Value = "$/unidad"
$ws.Range("O129").Value = "Región Metropolitana"
$ws.Range("P129").Value = 878
$ws.Range("Q129").Value = 1
$ws.Range("R129").Value = "Hortaliza"
